$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.036.20'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.834.15'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6271'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9996'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07607'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.69%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.90'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2927'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.59'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07722'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.838.54'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.953'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6651'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001016'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +17.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '82.83'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.056'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.026.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '226.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.192'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.56'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.503'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1372'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.37%  '
$ws.Range("E29").Value = '  -1.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.100'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.015'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.190'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05231'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.844'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7346'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.139'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.699'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.238.55'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.70%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.757'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01783'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.369'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8962'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9998'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.03'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000127'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.981.85'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '64.12'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5109'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4037'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.876'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05762'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.11%  '
